$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Num." index column (A) for rows 15 through 291: add 9 to each
# existing value (1..277 -> 10..286), removing the now-unused leading rows'
# worth of numbering so the sequence continues from the kept block above.
for ($r = 15; $r -le 291; $r++) {
    $cell = $ws.Range("A$r")
    $cell.Value = $cell.Value2 + 9
}

# Update the sheet view: clear the scrolled "topLeftCell" and move the
# active selection to A15.
$window = $excel.ActiveWindow
$window.ScrollRow = 1
$window.ScrollColumn = 1
$ws.Range("A15").Select()
